# Update "Estado de Cuenta" (NIT-9006476081): remove the two stale worker/period
# records for JAIRO ENRIQUE MELGAREJO ALVAREZ (periods 2507 and 2506), keeping
# only GERMAN VICENTE CASIANI HURTADO (period 2303), and refresh the summary
# totals (Valor Mora, Cant. Trabajadores, Cant. Periodos) to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 17 held the second JAIRO period (2506); row 16 held the first (2507).
# Deleting both shifts the GERMAN row (originally row 18) up to row 16, and
# shifts the trailing signature block up along with it.
$ws.Rows(17).Delete()
$ws.Rows(16).Delete()

# Refresh the rolled-up summary figures for the remaining single worker/period.
$ws.Range("E11").Value = 10827
$ws.Range("C13").Value = 1
$ws.Range("F13").Value = 1

# Column D narrows slightly now that the longest name left is shorter.
# (ColumnWidth round-trips with a fixed +5/6 offset in this host, so back it out
# here to land exactly on the target stored width of 34.)
$ws.Columns("D").ColumnWidth = 33.1666666666667
